$wb = $excel.ActiveWorkbook

# --- Sheet "levers": build out the expanded lever table (quartil steps) ---
$ws = $wb.Worksheets.Item("levers")
$ws.Activate()

# Row 2 (head of the two shared-formula runs)
$ws.Range("A2").Value = 1
$ws.Range("B2").Formula = "=""ADV-""&D2"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Formula = "=""ADV-""&D3"
$ws.Range("C3").Value = 1
$ws.Range("D3").Formula = "=D2+0.5"

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Formula = "=""ADV-""&D4"
$ws.Range("C4").Value = 1
$ws.Range("D4").Formula = "=D3+0.5"

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Formula = "=""ADV-""&D5"
$ws.Range("C5").Value = 1
$ws.Range("D5").Formula = "=D4+0.5"

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Formula = "=""ADV-""&D6"
$ws.Range("C6").Value = 1
$ws.Range("D6").Formula = "=D5+0.5"

# Row 7
$ws.Range("A7").Value = 6
$ws.Range("B7").Formula = "=""ADV-""&D7"
$ws.Range("C7").Value = 1
$ws.Range("D7").Formula = "=D6+0.5"

# Row 8
$ws.Range("A8").Value = 7
$ws.Range("B8").Formula = "=""ADV-""&D8"
$ws.Range("C8").Value = 1
$ws.Range("D8").Formula = "=D7+0.5"

# Row 9 (the original "ADV-0 / off" row, pushed down from row 3)
$ws.Range("A9").Value = 8
$ws.Range("B9").Formula = "=""ADV-""&D9"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0

$ws.Range("D10").Select()

# --- Sheet "params": stays the active/selected tab; just move the active cell to C11 ---
$wsParams = $wb.Worksheets.Item("params")
$wsParams.Activate()
$wsParams.Range("C11").Select()
